{"js": "// Change: \"...spoon that one uses to temper sand...\"\n//      -> \"...spoon with which one tempers the sand...\"\n// (paragraph <head>Cleaning the <tl>bowl</tl> and the <tl>spoon</tl> ... sand</head>)\n//\n// The surrounding text repeats a very similar phrase (\"...spoon with which\n// you temper your sand...\") in a different paragraph, so all lookups are\n// scoped to the specific paragraph that contains \"that one uses\" to avoid\n// touching the unrelated paragraph.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"that one uses\") !== -1) {\n    target = p;\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the target paragraph containing 'that one uses'.\");\n}\n\n// Step 1: the run \" \" (before \"that\") together with the runs \"that one uses\"\n// and \" to \" become a single run \" with which \" (keeps the leading run's\n// formatting: color 000000 / rtl 0).\nlet matches = target.search(\" that one uses to \", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\nif (matches.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for ' that one uses to ', found \" + matches.items.length);\n}\nmatches.items[0].insertText(\" with which \", Word.InsertLocation.replace);\nawait context.sync();\n\n// Step 2: the run \"temper\" becomes \"one tempers the\" (keeps its own\n// formatting: rtl 0, no explicit color).\nmatches = target.search(\"temper\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\nif (matches.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'temper', found \" + matches.items.length);\n}\nmatches.items[0].insertText(\"one tempers the\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change: \"...spoon that one uses to temper sand...\"\n#      -> \"...spoon with which one tempers the sand...\"\n# (paragraph <head>Cleaning the <tl>bowl</tl> and the <tl>spoon</tl> ... sand</head>)\n#\n# The surrounding document repeats a very similar phrase (\"...spoon with\n# which you temper your sand...\") in a different paragraph, so every Find\n# is scoped to the specific paragraph that contains \"that one uses\" to\n# avoid touching the unrelated paragraph.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*that one uses*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the target paragraph containing 'that one uses'.\"\n}\n\n# Step 1: the run \" \" (before \"that\") together with the runs \"that one uses\"\n# and \" to \" become a single run \" with which \" (keeps the leading run's\n# formatting: color 000000 / rtl 0).\n$r1 = $target.Range\n$r1.Find.ClearFormatting()\n$found1 = $r1.Find.Execute(\" that one uses to \", $false, $false, $false, $false, $false, $true, 1, $false, \" with which \", 2)\nif (-not $found1) {\n    throw \"Could not find ' that one uses to ' in the target paragraph.\"\n}\n\n# Step 2: the run \"temper\" becomes \"one tempers the\" (keeps its own\n# formatting: rtl 0, no explicit color). Re-resolve the paragraph since the\n# text length shifted after step 1.\n$target2 = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*with which*\") {\n        $target2 = $p\n        break\n    }\n}\nif ($target2 -eq $null) {\n    throw \"Could not find the paragraph after the first replacement.\"\n}\n\n$r2 = $target2.Range\n$r2.Find.ClearFormatting()\n$found2 = $r2.Find.Execute(\"temper\", $false, $false, $false, $false, $false, $true, 1, $false, \"one tempers the\", 2)\nif (-not $found2) {\n    throw \"Could not find 'temper' in the target paragraph.\"\n}\n"}
